# Modifs in database (2)
# Halve the capacity values in columns B and C (rows 2-32) of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 2).Value2 / 2
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 3).Value2 / 2
}
